# Generate Report for Handback
#
# The localization-status report is regenerated after a handback event for
# file "5ffd0142-ccb6-4c91-9c0e-c04c855f5fef.md": its status flips from
# "Ready for handoff" to "Handed back: in sync with en-US" on every sheet
# that tracks it, the per-locale "Latest Handback DateTime" is stamped with
# the new handback time, and the stale "handback file is not the latest"
# error note is cleared now that the handback is current.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row for 5ffd0142-ccb6-4c91-9c0e-c04c855f5fef.md ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"   # zh-cn status
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"   # de-de status

# --- zh-cn sheet: row for 5ffd0142-ccb6-4c91-9c0e-c04c855f5fef.md ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"       # Status
$wsZhCn.Range("K3").Value = "2016-08-18 06:46:53"                   # Latest Handback DateTime
$wsZhCn.Range("P3").Value = ""                                      # Error Detail cleared
$wsZhCn.Columns.Item(16).AutoFit()

# --- de-de sheet: row for 5ffd0142-ccb6-4c91-9c0e-c04c855f5fef.md ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"       # Status
$wsDeDe.Range("K3").Value = "2016-08-18 06:47:04"                   # Latest Handback DateTime
$wsDeDe.Range("P3").Value = ""                                      # Error Detail cleared
$wsDeDe.Columns.Item(16).AutoFit()
